$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2703
$ws1.Range("F9").Value = 1588
$ws1.Range("F10").Value = 7306
$ws1.Range("F12").Value = 7450
$ws1.Range("F14").Value = 29
$ws1.Range("F15").Value = 5867
$ws1.Range("F16").Value = 3184
$ws1.Range("F17").Value = 3557
$ws1.Range("F21").Value = 218
$ws1.Range("F22").Value = 2014
$ws1.Range("F23").Value = 94
$ws1.Range("F29").Value = 2526
$ws1.Range("F30").Value = 1353
$ws1.Range("F31").Value = 3052
$ws1.Range("F34").Value = 191
$ws1.Range("F37").Value = 216
$ws1.Range("F38").Value = 510

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 49
$ws2.Range("F7").Value = 38
$ws2.Range("F11").Value = 28
$ws2.Range("F14").Value = 7
$ws2.Range("F15").Value = 28
$ws2.Range("F16").Value = 56
$ws2.Range("F17").Value = 75

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 102

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 49
$ws4.Range("F11").Value = 102
$ws4.Range("F12").Value = 2703
$ws4.Range("F13").Value = 1588
$ws4.Range("F15").Value = 7306
$ws4.Range("F17").Value = 7450
$ws4.Range("F19").Value = 29
$ws4.Range("F20").Value = 5867
$ws4.Range("F21").Value = 3184
$ws4.Range("F22").Value = 3557
$ws4.Range("F25").Value = 28
$ws4.Range("F29").Value = 2014
$ws4.Range("F30").Value = 7
$ws4.Range("F31").Value = 28
$ws4.Range("F32").Value = 56
$ws4.Range("F38").Value = 2526
$ws4.Range("F39").Value = 1353
$ws4.Range("F40").Value = 75
$ws4.Range("F41").Value = 3052
$ws4.Range("F44").Value = 191
$ws4.Range("F48").Value = 510
